$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 199, shifting existing rows 199-201 down to 200-202
$ws.Rows.Item(199).Insert()

# Fill the newly inserted row 199 with the new weekly data
$ws.Range("A199").Value = 11
$ws.Range("B199").Value = "Vega Monumental Concepción"
$ws.Range("C199").Value = "Bíobío"
$ws.Range("D199").Value = 45239
$ws.Range("D199").NumberFormat = $ws.Range("D200").NumberFormat
$ws.Range("E199").Value = 8
$ws.Range("F199").Value = 100112028
$ws.Range("G199").Value = "Sandia"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 500
$ws.Range("K199").Value = 700
$ws.Range("L199").Value = 750
$ws.Range("M199").Value = 730
$ws.Range("N199").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O199").Value = "Perú"
$ws.Range("P199").Value = 730
$ws.Range("Q199").Value = 1
$ws.Range("R199").Value = "Hortaliza"
